# Group the existing "Group 12" shape together with the "TextBox 5" shape
# on slide 3 into a new enclosing group ("Group 1"). This mirrors selecting
# both shapes in the UI and pressing Ctrl+G: PowerPoint creates a brand new
# outer group shape, re-parenting the two original shapes underneath it,
# while leaving their own ids/names/contents untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$range = $s.Shapes.Range(@("Group 12", "TextBox 5"))
$range.Group() | Out-Null
